# Auto-generated edit script applying cryptos list value updates
# (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.700.93"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "2.494.97"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.26"
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.95"
$ws.Range("D4").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E6").Value = "  +3.21%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.516"
$ws.Range("D4").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.143"
$ws.Range("D4").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("E9").Value = "  +6.35%  "
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.341"
$ws.Range("D4").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("E11").Value = "  +4.10%  "
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").Value = "2.945.86"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.83"
$ws.Range("D4").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").Value = "  +2.98%  "
$ws.Range("D15").Value = "67.590.59"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("D4").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "  +3.04%  "
$ws.Range("D17").Value = "2.496.54"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.11"
$ws.Range("D4").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.52"
$ws.Range("D4").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E19").Value = "  +2.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "352.08"
$ws.Range("D4").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.63"
$ws.Range("D4").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = "  +3.25%  "
$ws.Range("E24").Value = "  +2.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.76"
$ws.Range("D4").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.20"
$ws.Range("D4").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("D27").Value = "2.628.57"
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D4").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "0.0₃0913"
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "509.45"
$ws.Range("D4").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.83"
$ws.Range("D4").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "  +3.45%  "
$ws.Range("E32").Value = "  +4.43%  "
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.124"
$ws.Range("D4").Copy()
$ws.Range("D35").PasteSpecial(-4122)
$ws.Range("E35").Value = "  +8.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.14"
$ws.Range("D4").Copy()
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("E36").Value = "  +2.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.49"
$ws.Range("D4").Copy()
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("E39").Value = "  +1.78%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D4").Copy()
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.75"
$ws.Range("D4").Copy()
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("E41").Value = "  +5.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.331"
$ws.Range("D4").Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("E42").Value = "  +2.18%  "
$ws.Range("E43").Value = "  +2.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.42"
$ws.Range("D4").Copy()
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("E44").Value = "  +3.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "146.04"
$ws.Range("D4").Copy()
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("E45").Value = "  +3.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.52"
$ws.Range("D4").Copy()
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("E46").Value = "  +3.02%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0259"
$ws.Range("E47").Value = "  +4.19%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.516"
$ws.Range("D4").Copy()
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("E49").Value = "  +2.82%  "
$ws.Range("E50").Value = "  +2.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.586"
$ws.Range("D4").Copy()
$ws.Range("D51").PasteSpecial(-4122)
$ws.Range("E51").Value = "  +1.32%  "
